$d = $word.ActiveDocument

$d.Content.Find.Execute("ابراهيم علي المتولى", $true, $false, $false, $false, $false,
                         $true, 1, $false, "نرمين رضا", 2)

$d.Content.Find.Execute("21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "100", 2)

$d.Content.Find.Execute("517", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1675", 2)

$d.Content.Find.Execute("04/11/2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "01/06/2021", 2)
